# Day5-Third commit: appointment booking / confirmation details.
#
# Adds a second booked-appointment confirmation email ("abi@.com") in D3,
# mirroring the existing confirmation email in D2 ("abi@gmail.com"):
#   - new shared string "abi@.com"
#   - D3 gets that value, styled with the built-in "Hyperlink" style
#   - a mailto: hyperlink is attached to D3
#   - the active selection moves to D9 (matching the recorded sheet view)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New appointment confirmation email value for row 3.
$ws.Range("D3").Value = "abi@.com"

# Hook up the mailto hyperlink (same pattern as the existing D2 entry).
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:abi@.com")

# Make sure D3 carries the same "Hyperlink" cell style as D2.
$ws.Range("D3").Style = "Hyperlink"

# Move/record the active selection at D9.
$ws.Range("D9").Select()
